$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.047.77'
$ws.Range('E2').Value = '  -0.06%  '

# Row 3
$ws.Range('D3').Value = '1.649.83'
$ws.Range('E3').Value = '  +0.48%  '

# Row 4
$ws.Range('E4').Value = '  -0.29%  '

# Row 5
$ws.Range('D5').Value = '218.15'
$ws.Range('E5').Value = '  +0.42%  '

# Row 6
$ws.Range('D6').Value = '0.5198'
$ws.Range('E6').Value = '  +0.41%  '

# Row 7
$ws.Range('E7').Value = '  -0.30%  '

# Row 8
$ws.Range('E8').Value = '  +0.96%  '

# Row 9
$ws.Range('D9').Value = '0.06323'
$ws.Range('E9').Value = '  +0.79%  '

# Row 10
$ws.Range('D10').Value = '20.38'
$ws.Range('E10').Value = '  +0.27%  '

# Row 11
$ws.Range('D11').Value = '0.07658'
$ws.Range('E11').Value = '  -1.26%  '

# Row 12
$ws.Range('D12').Value = '4.583'
$ws.Range('E12').Value = '  +2.54%  '

# Row 13
$ws.Range('D13').Value = '1.645.04'
$ws.Range('E13').Value = '  -1.84%  '

# Row 14
$ws.Range('D14').Value = '1.876.47'
$ws.Range('E14').Value = '  +0.43%  '

# Row 15
$ws.Range('D15').Value = '0.5584'
$ws.Range('E15').Value = '  +0.62%  '

# Row 16
$ws.Range('D16').Value = '0.0₅8125'
$ws.Range('E16').Value = '  +2.03%  '

# Row 17
$ws.Range('D17').Value = '65.21'
$ws.Range('E17').Value = '  +0.89%  '

# Row 18
$ws.Range('D18').Value = '26.032.14'
$ws.Range('E18').Value = '  -0.09%  '

# Row 19
$ws.Range('E19').Value = '  -0.24%  '

# Row 20
$ws.Range('D20').Value = '4.612'
$ws.Range('E20').Value = '  +0.12%  '

# Row 21
$ws.Range('D21').Value = '10.51'
$ws.Range('E21').Value = '  +4.47%  '

# Row 22
$ws.Range('D22').Value = '191.23'

# Row 23
$ws.Range('D23').Value = '5.900'
$ws.Range('E23').Value = '  -0.58%  '

# Row 24
$ws.Range('E24').Value = '  -0.26%  '

# Row 25
$ws.Range('D25').Value = '143.51'
$ws.Range('E25').Value = '  -2.20%  '

# Row 26
$ws.Range('D26').Value = '0.1183'
$ws.Range('E26').Value = '  -1.41%  '

# Row 27
$ws.Range('D27').Value = '7.185'
$ws.Range('E27').Value = '  +0.57%  '

# Row 28
$ws.Range('D28').Value = '15.85'
$ws.Range('E28').Value = '  -0.05%  '

# Row 29
$ws.Range('D29').Value = '1.507'
$ws.Range('E29').Value = '  +1.78%  '

# Row 30
$ws.Range('D30').Value = '0.05369'
$ws.Range('E30').Value = '  -4.73%  '

# Row 31
$ws.Range('D31').Value = '1.268'
$ws.Range('E31').Value = '  +0.18%  '

# Row 32
$ws.Range('D32').Value = '3.454'
$ws.Range('E32').Value = '  -0.24%  '

# Row 33
$ws.Range('D33').Value = '3.344'
$ws.Range('E33').Value = '  -0.39%  '

# Row 34
$ws.Range('E34').Value = '  -2.11%  '

# Row 35
$ws.Range('D35').Value = '2.421'

# Row 36
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').Value = '2.778'
$ws.Range('E36').Value = '  -0.34%  '

# Row 37
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = '0.9446'
$ws.Range('E37').Value = '  +1.06%  '

# Row 38
$ws.Range('D38').Value = '0.5631'
$ws.Range('E38').Value = '  -0.15%  '

# Row 39
$ws.Range('D39').Value = '0.01576'
$ws.Range('E39').Value = '  +0.19%  '

# Row 40
$ws.Range('D40').Value = '5.892'
$ws.Range('E40').Value = '  -0.69%  '

# Row 42
$ws.Range('D42').Value = '1.030.12'
$ws.Range('E42').Value = '  -2.49%  '

# Row 43
$ws.Range('D43').Value = '0.8257'
$ws.Range('E43').Value = '  -1.51%  '

# Row 44
$ws.Range('D44').Value = '100.66'
$ws.Range('E44').Value = '  -1.86%  '

# Row 45
$ws.Range('D45').Value = '1.786.29'
$ws.Range('E45').Value = '  +0.35%  '

# Row 46
$ws.Range('E46').Value = '  +5.48%  '

# Row 47
$ws.Range('D47').Value = '57.21'
$ws.Range('E47').Value = '  +0.62%  '

# Row 48
$ws.Range('D48').Value = '1.003'
$ws.Range('E48').Value = '  -0.40%  '

# Row 49
$ws.Range('E49').Value = '  -0.38%  '

# Row 50
$ws.Range('D50').Value = '7.906'
$ws.Range('E50').Value = '  -0.50%  '

# Row 51
$ws.Range('D51').Value = '0.05139'
$ws.Range('E51').Value = '  -3.68%  '
